$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (ID 6): mark the "Solved" flag on
$ws.Range("E8").Value = 1

# Row 10 (ID 8): new bug entry "Enemy Healthbar BUG"
$ws.Range("B10").Value = "Enemy Healthbar BUG"
$ws.Range("C10").Value = "Balken ist teilweise kurz wenn man draufschlägt, steigt dann rapide an und decreast auf den eigentlichen Wert"
$ws.Range("E10").Value = 1

# Row 12 (ID 10): new bug entry "Stoney_Dialoge"
$ws.Range("B12").Value = "Stoney_Dialoge"
$ws.Range("C12").Value = "Dialoge wird öfters gestartet"
$ws.Range("D12").Value = "vermutlich CoRoutine mehrfach gestartet"

# Update the active selection to match where the author left off editing
$ws.Range("B13").Select()
